$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Library_Formula")
$ws.Rows("139:139").Insert()
$ws.Range("A139").Value = "CREATE/MODIFY"
$ws.Range("B139").Value = "LIB_EWS_IT"
$ws.Range("C139").Value = "INDICATOR_88"
$ws.Range("E139").Value = "String"
$ws.Range("F139").Value = "String"

$ws.Activate()
$ws.Range("B137").Select() | Out-Null
